$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: update the intro text in A1 to drop the "phải học ĐHKHTN" (must be
# a ĐHKHTN student) requirement wording - no longer need KHTN students.
$ws.Range("A1").Value = "Bước 1: Điền vài thông tin của bạn, thông tin này để kiểm tra lại khi cần"

# Step 2: the survey no longer needs the MSSV (student id) field. Remove the
# whole "MSSV / 0812000" row (row 2); everything below shifts up by one row
# automatically (values, hyperlinks, and used styles renumber accordingly).
$ws.Rows("2:2").Delete()
